$d = $word.ActiveDocument

# Anchor on the paragraph containing "LOB1053: Física III (Requisito)" so the
# deletion is located relative to content rather than a fixed paragraph index.
$anchor = $d.Content.Duplicate
[void]$anchor.Find.Execute("LOB1053: Física III (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorStart = $anchor.Paragraphs.Item(1).Range.Start

$allParas = $d.Paragraphs
$n = $allParas.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $allParas.Item($i)
    if ($para.Range.Start -eq $anchorStart) {
        # The three paragraphs right after the anchor are removed in full:
        #   (blank "Normal" paragraph)
        #   "Ver no Jupiter Salvar em pdf Salvar em docx"
        #   "© 2020 . Contact: ... Creative Commons Attribution"
        $deleteStart = $allParas.Item($i + 1).Range.Start
        $deleteEnd = $allParas.Item($i + 4).Range.Start
        $r = $d.Range($deleteStart, $deleteEnd)
        $r.Delete()
        break
    }
}
